$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the generated credential values in row 2 (newly added iAuthor TC's)
$ws.Range("A2").Value = "ctrTv372"
$ws.Range("B2").Value = 231011251
$ws.Range("C2").Value = "xcdylvi50"
$ws.Range("D2").Value = "EFrk&48#"
$ws.Range("F2").Value = "pgUnTswa"
$ws.Range("G2").Value = "uzgK"
